$d = $word.ActiveDocument

$pairs = @(
    @("148×4=", "571×7="),
    @("361×5=", "914×3="),
    @("738×8=", "376×7="),
    @("183×6=", "396×3="),
    @("970×3=", "147×4="),
    @("497×5=", "209×4="),
    @("297×9=", "606×8="),
    @("828×3=", "456×9="),
    @("147×7=", "837×7="),
    @("279×9=", "932×3="),
    @("387×9=", "704×9="),
    @("244×8=", "526×6="),
    @("466×8=", "157×7="),
    @("640×6=", "300×7="),
    @("392×5=", "974×2="),
    @("755×3=", "325×6="),
    @("880×2=", "166×8="),
    @("182×7=", "661×5="),
    @("318×6=", "169×9="),
    @("294×6=", "129×5="),
    @("379×3=", "455×3="),
    @("770×6=", "401×8="),
    @("631×2=", "560×9="),
    @("670×8=", "326×6="),
    @("561×7=", "634×4=")
)

foreach ($pair in $pairs) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
